$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "Paid Date" / "Outstanding" columns right by one.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour
# (column M, "In Advance") the way Excel does when inserting a column.
$ws.Columns("N").ColumnWidth = 10.333333333333334

# Make "Repayment Schedule" the active sheet/tab and select cell P6 on it
# (previously "Transactions" was the active tab).
$ws.Activate() | Out-Null
$ws.Range("P6").Select() | Out-Null
